$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

# Delete entire row 24 (shifts all rows below it up by one),
# which removes the old "6141 EL PAMPERO 2618" record and
# causes every subsequent record to move up one row, dropping
# the last row (old row 73) off the bottom of the table.
$ws.Rows.Item(24).Delete()
